$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Wed Dec 20 12:56:49 EST 2023"
$ws.Range("B3").Value = "Wed Dec 20 12:57:02 EST 2023"
$ws.Range("B5").Value = "Wed Dec 20 12:57:14 EST 2023"
